$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I6").Value = "sv"
$ws.Range("J6").Value = "Statement-opinion"

# Row 41: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I41").Value = "sv"
$ws.Range("J41").Value = "Statement-opinion"

# Row 80: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I80").Value = "sv"
$ws.Range("J80").Value = "Statement-opinion"

# Row 102: b/Acknowledge (Backchannel) -> aa/Agree/Accept
$ws.Range("I102").Value = "aa"
$ws.Range("J102").Value = "Agree/Accept"

# Row 124: sd/Statement-non-opinion -> aa/Agree/Accept
$ws.Range("I124").Value = "aa"
$ws.Range("J124").Value = "Agree/Accept"

# Row 126: sd/Statement-non-opinion -> aa/Agree/Accept
$ws.Range("I126").Value = "aa"
$ws.Range("J126").Value = "Agree/Accept"

# Row 133: sd/Statement-non-opinion -> aa/Agree/Accept
$ws.Range("I133").Value = "aa"
$ws.Range("J133").Value = "Agree/Accept"

# Row 134: sv/Statement-opinion -> sd/Statement-non-opinion
$ws.Range("I134").Value = "sd"
$ws.Range("J134").Value = "Statement-non-opinion"
